$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Update the problem text box: new wording + blank spacer lines + bonus question + hint ---
$contentShape = $s.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$tr.Text = "How many squares are there in the following 3x3 grid?" + "`r`r`r`r" + `
    "Bonus question: how many squares are there in a 8x8 grid?" + "`r" + `
    "Hint: look for number patterns in 3x3 grid."

# Indent the hint line (6th paragraph) to the second outline level.
$tr.Paragraphs(6, 1).IndentLevel = 2

# --- Shrink/reposition the grid diagram group to make room for the extra lines of text ---
# (Left/Top/Width/Height are expressed in points; the literals below are chosen so that the
# engine's internal float32 storage round-trips to the exact target EMU values.)
$group = $s.Shapes.Item(3)
$group.Left = 228.2263779527559
$group.Top = 243.50945291889764
$group.Width = 86.60377952755906
$group.Height = 86.60377952755906
